$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.393.70"
$ws.Range("E2").Value = "  -3.06%  "
$ws.Range("D3").Value = "3.317.71"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'558.55"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.23%  "
$ws.Range("D6").Value = "'142.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.96%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.319.96"
$ws.Range("E8").Value = "  -3.08%  "
$ws.Range("D9").Value = "'0.476"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("D10").Value = "'7.88"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("E11").Value = "  -3.60%  "
$ws.Range("D12").Value = "'0.409"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "'27.06"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.54%  "
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("D18").Value = "60.358.17"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").Value = "'6.20"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("D20").Value = "'14.50"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "'8.67"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "'376.57"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("D23").Value = "'74.15"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("E24").Value = "  -3.76%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "3.436.40"
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("E27").Value = "  -6.63%  "
$ws.Range("E28").Value = "  -4.99%  "
$ws.Range("D30").Value = "'7.32"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.12%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("D34").Value = "'22.61"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("E35").Value = "  -3.82%  "
$ws.Range("D36").Value = "'5.19"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.97%  "
$ws.Range("D37").Value = "'1.54"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.58%  "
$ws.Range("D38").Value = "'166.73"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").Value = "'6.76"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("E41").Value = "  -14.49%  "
$ws.Range("D42").Value = "'0.0741"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.62%  "
$ws.Range("D43").Value = "'42.03"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("E45").Value = "  -3.67%  "
$ws.Range("E46").Value = "  -4.80%  "
$ws.Range("D47").Value = "'1.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("D48").Value = "2.368.37"
$ws.Range("E48").Value = "  -6.90%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "'6.55"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("D51").Value = "'21.42"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.79%  "
